$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AgTests (H) and AgPosit (I) values for rows 286-321
$ws.Range("H286").Value = 54215
$ws.Range("I286").Value = 4250
$ws.Range("H287").Value = 57631
$ws.Range("I287").Value = 3917
$ws.Range("H288").Value = 56088
$ws.Range("I288").Value = 3959
$ws.Range("H289").Value = 64387
$ws.Range("I289").Value = 3698
$ws.Range("H292").Value = 81217
$ws.Range("I292").Value = 7191
$ws.Range("H293").Value = 81670
$ws.Range("I293").Value = 5794
$ws.Range("H294").Value = 90791
$ws.Range("I294").Value = 5028
$ws.Range("H297").Value = 2282
$ws.Range("I297").Value = 208
$ws.Range("H299").Value = 63795
$ws.Range("I299").Value = 6698
$ws.Range("H300").Value = 70477
$ws.Range("I300").Value = 6917
$ws.Range("H301").Value = 69571
$ws.Range("I301").Value = 5560
$ws.Range("H302").Value = 72559
$ws.Range("I302").Value = 5284
$ws.Range("H306").Value = 70727
$ws.Range("I306").Value = 7182
$ws.Range("H307").Value = 73350
$ws.Range("I307").Value = 6323
$ws.Range("H309").Value = 57091
$ws.Range("I309").Value = 3950
$ws.Range("H310").Value = 90811
$ws.Range("I310").Value = 5414
$ws.Range("H313").Value = 72895
$ws.Range("I313").Value = 3546
$ws.Range("H314").Value = 65072
$ws.Range("I314").Value = 3346
$ws.Range("H315").Value = 65997
$ws.Range("I315").Value = 3106
$ws.Range("H316").Value = 49041
$ws.Range("H317").Value = 61020
$ws.Range("H320").Value = 83674
$ws.Range("I320").Value = 4065
$ws.Range("H321").Value = 80975
$ws.Range("I321").Value = 2654

# Add new row 322 for 2021-01-20 data
$ws.Range("A322").Value = 44216
$ws.Range("B322").Value = 231242
$ws.Range("C322").Value = 181129
$ws.Range("D322").Value = 46312
$ws.Range("E322").Value = 12668
$ws.Range("F322").Value = 2464
$ws.Range("G322").Value = 3801
$ws.Range("H322").Value = 82510
$ws.Range("I322").Value = 2045
